# Apply the weekly cryptos data refresh (GitHub Actions scrape) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.653.70"
$ws.Range("E2").Value = "  +2.39%  "

# Row 3
$ws.Range("D3").Value = "3.027.95"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D5").Value = "510.28"
$ws.Range("E5").Value = "  +2.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D6").Value = "139.31"
$ws.Range("E6").Value = "  +4.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D8").Value = "0.441"
$ws.Range("E8").Value = "  +3.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D9").Value = "7.50"
$ws.Range("E9").Value = "  +2.12%  "

# Row 10
$ws.Range("E10").Value = "  +3.90%  "

# Row 11
$ws.Range("E11").Value = "  +5.12%  "

# Row 12
$ws.Range("D12").Value = "3.548.22"
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("E13").Value = "  +2.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D14").Value = "26.61"
$ws.Range("E14").Value = "  +5.76%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").Value = "  +10.52%  "

# Row 16
$ws.Range("D16").Value = "57.601.38"
$ws.Range("E16").Value = "  +2.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D17").Value = "6.22"
$ws.Range("E17").Value = "  +10.01%  "

# Row 18
$ws.Range("D18").Value = "3.028.60"
$ws.Range("E18").Value = "  +2.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D19").Value = "12.90"
$ws.Range("E19").Value = "  +5.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D20").Value = "8.00"
$ws.Range("E20").Value = "  +4.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D21").Value = "333.21"

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D23").Value = "5.76"
$ws.Range("E23").Value = "  +0.93%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D24").Value = "0.498"
$ws.Range("E24").Value = "  +6.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D25").Value = "64.49"
$ws.Range("E25").Value = "  +4.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +4.47%  "

# Row 27
$ws.Range("E27").Value = "  -0.28%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  +5.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D29").Value = "6.81"
$ws.Range("E29").Value = "  +7.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +11.16%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D31").Value = "1.22"
$ws.Range("E31").Value = "  +3.96%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D32").Value = "1.80"
$ws.Range("E32").Value = "  +4.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D33").Value = "20.80"
$ws.Range("E33").Value = "  +3.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D34").Value = "155.55"
$ws.Range("E34").Value = "  -1.75%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D35").Value = "4.71"
$ws.Range("E35").Value = "  +7.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D36").Value = "5.86"
$ws.Range("E36").Value = "  +6.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  +2.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D38").Value = "24.61"
$ws.Range("E38").Value = "  +8.22%  "

# Row 39
$ws.Range("E39").Value = "  +2.46%  "

# Row 40
$ws.Range("D40").Value = "3.059.24"
$ws.Range("E40").Value = "  +1.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D41").Value = "37.33"
$ws.Range("E41").Value = "  +3.46%  "

# Row 42
$ws.Range("E42").Value = "  +9.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.16%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.306.03"
$ws.Range("E44").Value = "  +3.57%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D45").Value = "0.653"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("E46").Value = "  +3.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D47").Value = "0.988"
$ws.Range("E47").Value = "  +1.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D48").Value = "6.02"
$ws.Range("E48").Value = "  +5.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +3.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D50").Value = "19.64"
$ws.Range("E50").Value = "  +4.71%  "

# Row 51
$ws.Range("E51").Value = "  -3.80%  "

